# Add data for 2021-11-06
# - Rename the sheet / update the "through October NN" header from 28 -> 29
# - Bump a handful of monthly carjacking counts by neighborhood

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Sheet title and running header both track the reporting cut-off date.
$ws.Name = "Through 2021-10-29"
$ws.Range("B1").Value = "October 2021 (through October 29)"

# Updated counts (row = neighborhood, column = month)
$ws.Range("AZ2").Value = 6
$ws.Range("B3").Value = 14
$ws.Range("AZ3").Value = 4
$ws.Range("B6").Value = 6
$ws.Range("D11").Value = 8
$ws.Range("B12").Value = 3
$ws.Range("AP15").Value = 3
$ws.Range("BJ18").Value = 1
$ws.Range("V20").Value = 2
$ws.Range("AF25").Value = 2
$ws.Range("AZ27").Value = 1
$ws.Range("L38").Value = 8
$ws.Range("AZ38").Value = 5
$ws.Range("AZ40").Value = 2
$ws.Range("B42").Value = 3
$ws.Range("AP45").Value = 1
$ws.Range("L49").Value = 1
$ws.Range("B57").Value = 4
$ws.Range("B62").Value = 1
$ws.Range("L64").Value = 3
$ws.Range("AP67").Value = 3
$ws.Range("AZ72").Value = 1
$ws.Range("B80").Value = 4
